# JS-SPA-Self-Evaluation-Protocol.xlsx edit script
# Commit message: "delete work, add edit profile"
#   - fills in the student's profile info (username / name / GitHub link)
#   - fills in the scoring numbers for the "Basic Options" section
#   - fills in the two "infinite" rows' comment counts
#   - recalculates the Total Score formula

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Student info (SoftUni Student Info block) ---
$ws.Range("C4").Value = "AntoniyaIvanova"
$ws.Range("C5").Value = "Antoniya Ivanova"

# --- GitHub profile link (turns into a real hyperlink, like pasting a URL in Excel) ---
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/toniivanova/AngularProject")

# --- Days/Number of commits comments ---
$ws.Range("E8").Value = 16
$ws.Range("E9").Value = 24

# --- Basic Options (up to 130) scores ---
$ws.Range("C11").Value = 10
$ws.Range("C12").Value = 10
$ws.Range("C13").Value = 5
$ws.Range("C14").Value = 5
$ws.Range("C15").Value = 5
$ws.Range("C16").Value = 5
$ws.Range("C17").Value = 10
$ws.Range("C18").Value = 10
$ws.Range("C19").Value = 10
$ws.Range("C20").Value = 5
$ws.Range("C21").Value = 5
$ws.Range("C22").Value = 5
$ws.Range("C23").Value = 5
$ws.Range("C24").Value = 5
$ws.Range("C28").Value = 5
$ws.Range("C29").Value = 3
$ws.Range("C30").Value = 3
$ws.Range("C31").Value = 5
$ws.Range("C32").Value = 5

# Reselect like the saved file shows (cursor left on C30)
$ws.Range("C30").Select()

$wb.Save()
